$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.557.12"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "1.954.29"
$ws.Range("E3").Value = "  +0.57%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'243.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("E6").Value = "  +0.40%  "
$ws.Range("D7").Value = "'60.67"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +5.54%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "'0.376"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.03%  "
$ws.Range("D10").Value = "'0.0788"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -6.69%  "
$ws.Range("E11").Value = "  +0.16%  "
$ws.Range("D12").Value = "'14.24"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.81%  "
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").Value = "'0.830"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.24%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "2.243.86"
$ws.Range("E14").Value = "  +0.70%  "
$ws.Range("D15").Value = "'21.61"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.46%  "
$ws.Range("E16").Value = "  +1.24%  "
$ws.Range("D17").Value = "1.948.00"
$ws.Range("E17").Value = "  +0.17%  "
$ws.Range("D18").Value = "36.428.64"
$ws.Range("E18").Value = "  +0.05%  "
$ws.Range("D19").Value = "'69.59"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.34%  "
$ws.Range("E20").Value = "  -1.47%  "
$ws.Range("D21").Value = "'229.23"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.28%  "
$ws.Range("D22").Value = "'5.07"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.42%  "
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("E24").Value = "  +3.18%  "
$ws.Range("E25").Value = "  +1.73%  "
$ws.Range("D26").Value = "'0.143"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.81%  "
$ws.Range("D27").Value = "'9.19"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.13%  "
$ws.Range("D28").Value = "'160.07"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.97%  "
$ws.Range("D29").Value = "'19.28"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.41%  "
$ws.Range("D30").Value = "'1.31"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +18.47%  "
$ws.Range("E31").Value = "  +0.78%  "
$ws.Range("D33").Value = "'0.0611"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.88%  "
$ws.Range("D34").Value = "'4.46"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.41%  "
$ws.Range("E35").Value = "  +9.48%  "
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("E37").Value = "  +3.92%  "
$ws.Range("E38").Value = "  -1.06%  "
$ws.Range("D39").Value = "'5.46"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -12.31%  "
$ws.Range("B40").Value = "HuobiToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D40").Value = "'2.91"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.73%  "
$ws.Range("B41").Value = "Cronos"
$ws.Range("C41").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D41").Value = "'0.0955"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.61%  "
$ws.Range("D42").Value = "'1.17"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.60%  "
$ws.Range("E43").Value = "  +0.20%  "
$ws.Range("D44").Value = "'15.80"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.46%  "
$ws.Range("D45").Value = "1.361.06"
$ws.Range("E45").Value = "  +1.25%  "
$ws.Range("D46").Value = "'88.55"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.17%  "
$ws.Range("E47").Value = "  -0.25%  "
$ws.Range("E48").Value = "  -1.24%  "
$ws.Range("E49").Value = "  +0.57%  "
$ws.Range("D50").Value = "'45.88"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.03%  "
$ws.Range("D51").Value = "'3.53"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +10.65%  "
